$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2026/2027 forecast ("f") column headers to projection ("p")
$ws.Range("G1").Value = "2026p"
$ws.Range("H1").Value = "2027p"
